$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-23 Saturday", "2024-11-24 Sunday"),
    @("82×74=6068", "57×79=4503"),
    @("43×38=1634", "73×50=3650"),
    @("54×34=1836", "42×81=3402"),
    @("86×93=7998", "85×12=1020"),
    @("85×50=4250", "24×24=576"),
    @("35×50=1750", "81×61=4941"),
    @("86×14=1204", "24×46=1104"),
    @("95×55=5225", "48×34=1632"),
    @("25×67=1675", "38×29=1102"),
    @("69×35=2415", "57×87=4959"),
    @("12×39=468", "22×72=1584"),
    @("28×90=2520", "60×51=3060"),
    @("70×49=3430", "70×79=5530"),
    @("15×54=810", "53×24=1272"),
    @("21×56=1176", "55×73=4015"),
    @("18×92=1656", "80×51=4080"),
    @("25×96=2400", "77×67=5159"),
    @("31×30=930", "61×95=5795"),
    @("47×62=2914", "72×11=792"),
    @("84×65=5460", "57×18=1026"),
    @("76×82=6232", "56×55=3080"),
    @("33×55=1815", "28×86=2408"),
    @("26×19=494", "92×69=6348"),
    @("41×38=1558", "41×22=902"),
    @("48×20=960", "99×92=9108")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
